# The "Förändrad" (Changed) column C tracks the date a row was last
# refreshed by the automatic updater. This run bumps that date by one
# day (2023-09-09 -> 2023-09-10, serials 45178 -> 45179) for every
# existing data row (rows 2-236), leaving all other data untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C236").Value = 45179
